# Trade #28 closed at 2026-02-17 12:38:00 - unknown UNKNOWN +0.000%
#
# Updates the Summary, Strategy Status, All Trades and MarketMaking sheets
# to reflect the newly closed trade (#28 / zero-indexed 27 on MarketMaking).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.66   # Current Capital
$summary.Range("B4").Value = 0.65      # Total P&L $
$summary.Range("B5").Value = 0.46      # Total P&L %
$summary.Range("B6").Value = 28        # Total Trades
$summary.Range("B8").Value = 9         # Losing Trades
$summary.Range("B9").Value = 39.29     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.66     # Capital
$status.Range("D4").Value = 28         # Trades
$status.Range("E4").Value = 0.65       # P&L $
$status.Range("F4").Value = 0.66       # P&L %
$status.Range("G4").Value = 39.29      # Win Rate %

# ---------------------------------------------------------------------
# New trade row, appended to both "All Trades" and "MarketMaking" sheets
# ---------------------------------------------------------------------
$newTradeRow = @(28, "2026-02-17", "12:37:53", "MarketMaking", "UP", 0.66, 0.65, "CLOSED", -1.5152, -0.01, 100.66, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.14)
# Columns that hold date/time-looking text and must stay plain text instead
# of being auto-converted to Excel date/time serial numbers.
$textColumns = @(2, 3)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 29
    for ($i = 0; $i -lt $newTradeRow.Length; $i++) {
        $col = $i + 1
        $cell = $ws.Cells.Item($row, $col)
        if ($textColumns -contains $col) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $newTradeRow[$i]
    }
}
